$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.519.46"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "1.812.12"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("E4").Value = "  +0.80%  "
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").Value = "308.41"
$ws.Range("D7").Value = "0.4565"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").Value = "0.3663"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "0.07133"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "0.8795"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "0.07749"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "19.36"
$ws.Range("E12").Value = "  -3.52%  "
$ws.Range("D13").Value = "1.810.19"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "5.290"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "6.373"
$ws.Range("E15").Value = "  -2.15%  "
$ws.Range("D16").Value = "86.59"
$ws.Range("E16").Value = "  -5.11%  "
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "0.000008589"
$ws.Range("E18").Value = "  -3.66%  "
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "26.592.05"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").Value = "5.011"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "1.987"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Value = "151.50"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").Value = "2.057"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "112.86"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").Value = "4.844"
$ws.Range("E29").Value = "  -3.79%  "
$ws.Range("D30").Value = "0.08689"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").Value = "3.057"
$ws.Range("E31").Value = "  -2.67%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").Value = "0.7329"
$ws.Range("E33").Value = "  -4.30%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "2.680"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.118"
$ws.Range("E35").Value = "  -4.02%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").Value = "1.084"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("D38").Value = "0.01952"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").Value = "0.05114"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").Value = "2.893"
$ws.Range("D41").Value = "6.976"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Value = "0.4990"
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").Value = "0.1560"
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("D44").Value = "8.161"
$ws.Range("E44").Value = "  -3.80%  "
$ws.Range("D45").Value = "1.008"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").Value = "0.4598"
$ws.Range("E46").Value = "  -4.04%  "
$ws.Range("D47").Value = "10.01"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").Value = "101.05"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").Value = "0.06003"
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").Value = "64.38"
$ws.Range("E51").Value = "  -1.43%  "
